# Applies the "Batterywise analysis" edit described in the commit diff:
#  - Swap Starting/Ending SoC (%) values (B6/B7)
#  - Append unit suffixes to several row labels (column A)
#  - Swap the Highest/Lowest Cell Voltage labels+values (rows 16/17)
#  - Swap the lowest/highest cell temp labels (rows 28/29)
#  - Fill in the previously-blank "Difference in Temperature" value
#  - Remove old "Maximum BMS Temperature in C" row, shifting rows 32-42
#    up by one (row 31 onward), relabeling them and updating their
#    values, and appending a brand-new row 43
#    ("Time spent in 80-90 km/h").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 6 & 7: Starting/Ending SoC (%) values swap ---
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 24

# --- Rows 8-30: label text updates (values unchanged) ---
$ws.Range("A8").Value  = "Total distance covered (km)"
$ws.Range("A9").Value  = "Total energy consumption(WH/KM)"
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0

# Rows 16/17: Highest/Lowest Cell Voltage swap (label AND value swap)
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.336
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.013

$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"

$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 12

$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# Rows 28/29: lowest/highest cell temp label swap (values stay put)
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Rows 31-43: old row 31 ("Maximum BMS Temperature in C") is
#     dropped and everything below it shifts up by one row; write the
#     resulting labels/values directly (covers the shift + the new
#     row 43 appended at the end) ---
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.497764468055556

$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001567027064297505

$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 4.459798994974874

$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 4.070351758793969

$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 5.397822445561139

$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 10.73701842546064

$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 19.12060301507538

$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 27.34924623115578

$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 20.05025125628141

$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 7.24036850921273

$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 1.461474036850921

# --- New row 43 ---
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
